$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns retain text formatting so numeric-looking strings
# (prices, percentages) are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "42.692.91"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3
$ws.Range("D3").Value = "2.265.20"
$ws.Range("E3").Value = "  -0.46%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").Value = "249.74"
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("D6").Value = "0.640"
$ws.Range("E6").Value = "  +2.44%  "

# Row 7
$ws.Range("D7").Value = "77.23"
$ws.Range("E7").Value = "  +7.19%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "0.641"
$ws.Range("E9").Value = "  -3.14%  "

# Row 10
$ws.Range("D10").Value = "40.09"
$ws.Range("E10").Value = "  +3.18%  "

# Row 11
$ws.Range("E11").Value = "  +0.11%  "

# Row 12
$ws.Range("D12").Value = "7.33"
$ws.Range("E12").Value = "  -1.36%  "

# Row 13
$ws.Range("E13").Value = "  +0.80%  "

# Row 14
$ws.Range("D14").Value = "2.603.06"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15
$ws.Range("D15").Value = "14.97"
$ws.Range("E15").Value = "  +1.19%  "

# Row 16
$ws.Range("E16").Value = "  -2.63%  "

# Row 17
$ws.Range("D17").Value = "2.275.04"
$ws.Range("E17").Value = "  +0.26%  "

# Row 18
$ws.Range("D18").Value = "42.569.38"
$ws.Range("E18").Value = "  -0.44%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0992"
$ws.Range("E19").Value = "  -1.19%  "

# Row 20
$ws.Range("E20").Value = "  -2.21%  "

# Row 21
$ws.Range("D21").Value = "72.09"
$ws.Range("E21").Value = "  -1.45%  "

# Row 22
$ws.Range("D22").Value = "235.65"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("E23").Value = "  -0.90%  "

# Row 24
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("E25").Value = "  -7.08%  "

# Row 26
$ws.Range("E26").Value = "  -0.54%  "

# Row 27
$ws.Range("D27").Value = "2.37"
$ws.Range("E27").Value = "  -2.65%  "

# Row 28
$ws.Range("E28").Value = "  +2.42%  "

# Row 29
$ws.Range("D29").Value = "167.54"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("D30").Value = "20.93"
$ws.Range("E30").Value = "  -0.35%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "6.40"
$ws.Range("E31").Value = "  -1.42%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0858"
$ws.Range("E32").Value = "  +6.44%  "

# Row 33
$ws.Range("E33").Value = "  -3.34%  "

# Row 34
$ws.Range("D34").Value = "31.30"
$ws.Range("E34").Value = "  -1.82%  "

# Row 35
$ws.Range("E35").Value = "  +1.07%  "

# Row 36
$ws.Range("D36").Value = "4.55"
$ws.Range("E36").Value = "  +1.44%  "

# Row 37
$ws.Range("D37").Value = "4.72"
$ws.Range("E37").Value = "  -1.03%  "

# Row 38
$ws.Range("E38").Value = "  -2.98%  "

# Row 39
$ws.Range("D39").Value = "13.79"
$ws.Range("E39").Value = "  +7.38%  "

# Row 40
$ws.Range("E40").Value = "  -3.30%  "

# Row 41
$ws.Range("D41").Value = "5.87"
$ws.Range("E41").Value = "  +0.95%  "

# Row 42
$ws.Range("D42").Value = "0.208"
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("D43").Value = "61.31"
$ws.Range("E43").Value = "  -1.42%  "

# Row 44
$ws.Range("D44").Value = "108.25"
$ws.Range("E44").Value = "  +13.38%  "

# Row 45
$ws.Range("D45").Value = "8.86"
$ws.Range("E45").Value = "  -4.60%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.100"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "4.62"
$ws.Range("E47").Value = "  -9.00%  "

# Row 48
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.29%  "

# Row 49
$ws.Range("E49").Value = "  -2.31%  "

# Row 50
$ws.Range("B50").Value = "Bonk"
$ws.Range("C50").Value = "https://coinranking.com/coin/jCd_nuYCH+bonk-bonk"
$ws.Range("D50").Value = "0.0000344"
$ws.Range("E50").Value = "  +136.28%  "

# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "1.17"
$ws.Range("E51").Value = "  -2.47%  "
